$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 1013.7143
$ws.Cells.Item(4, 9).Value = 938.4
$ws.Cells.Item(4, 11).Value = 938.4
$ws.Cells.Item(4, 13).Value = -824.4
$ws.Cells.Item(29, 8).Value = 787.7143
$ws.Cells.Item(29, 10).Value = 4
$ws.Cells.Item(29, 12).Value = 12
$ws.Cells.Item(29, 14).Value = -574
$ws.Cells.Item(38, 8).Value = 769953.3
$ws.Cells.Item(38, 9).Value = 1250299.1
$ws.Cells.Item(38, 11).Value = 3750897.3
$ws.Cells.Item(38, 13).Value = -3750525.3
$ws.Cells.Item(40, 8).Value = 2170.25
$ws.Cells.Item(40, 9).Value = 1341.5
$ws.Cells.Item(40, 11).Value = 1341.5
$ws.Cells.Item(40, 13).Value = -1166.5
$ws.Cells.Item(53, 8).Value = 226.44444
$ws.Cells.Item(53, 9).Value = 286.2857
$ws.Cells.Item(53, 10).Value = 17
$ws.Cells.Item(53, 11).Value = 286.2857
$ws.Cells.Item(53, 12).Value = 17
$ws.Cells.Item(53, 13).Value = 350.7143
$ws.Cells.Item(53, 14).Value = -1291
$ws.Cells.Item(70, 8).Value = 3567.0833
$ws.Cells.Item(70, 9).Value = 7550.5
$ws.Cells.Item(70, 10).Value = 1575.375
$ws.Cells.Item(70, 11).Value = 22651.5
$ws.Cells.Item(70, 12).Value = 4726.125
$ws.Cells.Item(70, 13).Value = -22381.5
$ws.Cells.Item(70, 14).Value = -5266.125
$ws.Cells.Item(73, 8).Value = 3567.0833
$ws.Cells.Item(73, 9).Value = 7550.5
$ws.Cells.Item(73, 10).Value = 1575.375
$ws.Cells.Item(73, 11).Value = 22651.5
$ws.Cells.Item(73, 12).Value = 4726.125
$ws.Cells.Item(73, 13).Value = -21715.5
$ws.Cells.Item(73, 14).Value = -6598.125
$ws.Cells.Item(88, 8).Value = 2787.5833
$ws.Cells.Item(88, 9).Value = 2472.4
$ws.Cells.Item(88, 10).Value = 3012.7144
$ws.Cells.Item(88, 11).Value = 2472.4
$ws.Cells.Item(88, 12).Value = 3012.7144
$ws.Cells.Item(88, 13).Value = -2066.4
$ws.Cells.Item(88, 14).Value = -3824.7144
$ws.Cells.Item(91, 8).Value = 2787.5833
$ws.Cells.Item(91, 9).Value = 2472.4
$ws.Cells.Item(91, 10).Value = 3012.7144
$ws.Cells.Item(91, 11).Value = 2472.4
$ws.Cells.Item(91, 12).Value = 3012.7144
$ws.Cells.Item(91, 13).Value = -1068.4
$ws.Cells.Item(91, 14).Value = -5820.7144
$ws.Cells.Item(98, 8).Value = 1829.3636
$ws.Cells.Item(98, 9).Value = 879.125
$ws.Cells.Item(98, 10).Value = 4363.3335
$ws.Cells.Item(98, 11).Value = 879.125
$ws.Cells.Item(98, 12).Value = 4363.3335
$ws.Cells.Item(98, 13).Value = 618.875
$ws.Cells.Item(98, 14).Value = -7359.3335
$ws.Cells.Item(116, 8).Value = 5038.6665
$ws.Cells.Item(116, 9).Value = 4996.4
$ws.Cells.Item(116, 11).Value = 4996.4
$ws.Cells.Item(116, 13).Value = -1554.4
$ws.Cells.Item(122, 8).Value = 1829.3636
$ws.Cells.Item(122, 9).Value = 879.125
$ws.Cells.Item(122, 10).Value = 4363.3335
$ws.Cells.Item(122, 11).Value = 2637.375
$ws.Cells.Item(122, 12).Value = 13090.0005
$ws.Cells.Item(122, 13).Value = -187.375
$ws.Cells.Item(122, 14).Value = -17990.0005
$ws.Cells.Item(129, 8).Value = 573
$ws.Cells.Item(129, 9).Value = 573
$ws.Cells.Item(129, 11).Value = 1719
$ws.Cells.Item(129, 13).Value = 3281
$ws.Cells.Item(132, 8).Value = 5808.75
$ws.Cells.Item(132, 9).Value = 5808.75
$ws.Cells.Item(132, 11).Value = 17426.25
$ws.Cells.Item(132, 13).Value = -14896.25
$ws.Cells.Item(141, 8).Value = 14047.5
$ws.Cells.Item(141, 9).Value = 18095
$ws.Cells.Item(141, 10).Value = 10000
$ws.Cells.Item(141, 11).Value = 54285
$ws.Cells.Item(141, 12).Value = 30000
$ws.Cells.Item(141, 13).Value = -49105
$ws.Cells.Item(141, 14).Value = -40360
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9074.064
$ws.Cells.Item(32, 9).Value = 8530.896000000001
$ws.Cells.Item(32, 11).Value = 8530.896000000001
$ws.Cells.Item(32, 13).Value = -8243.896000000001
$ws.Cells.Item(61, 8).Value = 2011.125
$ws.Cells.Item(61, 9).Value = 2011.125
$ws.Cells.Item(61, 11).Value = 2011.125
$ws.Cells.Item(61, 13).Value = -1799.125
$ws.Cells.Item(102, 8).Value = 2349.6667
$ws.Cells.Item(102, 9).Value = 2349.6667
$ws.Cells.Item(102, 11).Value = 2349.6667
$ws.Cells.Item(102, 13).Value = -727.6667000000002
$ws.Cells.Item(132, 8).Value = 6998.2856
$ws.Cells.Item(132, 9).Value = 4989
$ws.Cells.Item(132, 10).Value = 7333.1665
$ws.Cells.Item(132, 11).Value = 14967
$ws.Cells.Item(132, 12).Value = 21999.4995
$ws.Cells.Item(132, 13).Value = -12437
$ws.Cells.Item(132, 14).Value = -27059.4995
$ws.Cells.Item(136, 8).Value = 2011.125
$ws.Cells.Item(136, 9).Value = 2011.125
$ws.Cells.Item(136, 11).Value = 6033.375
$ws.Cells.Item(136, 13).Value = -3483.375
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2883.15
$ws.Cells.Item(86, 9).Value = 2946.1765
$ws.Cells.Item(86, 11).Value = 2946.1765
$ws.Cells.Item(86, 13).Value = -1823.1765
$ws.Cells.Item(89, 8).Value = 2883.15
$ws.Cells.Item(89, 9).Value = 2946.1765
$ws.Cells.Item(89, 11).Value = 14730.8825
$ws.Cells.Item(89, 13).Value = -9114.8825
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 1583
$ws.Cells.Item(4, 9).Value = 1749
$ws.Cells.Item(4, 11).Value = 1749
$ws.Cells.Item(4, 13).Value = -1637
$ws.Cells.Item(31, 8).Value = 2522.3333
$ws.Cells.Item(31, 9).Value = 2462.625
$ws.Cells.Item(31, 11).Value = 2462.625
$ws.Cells.Item(31, 13).Value = -2167.625
$ws.Cells.Item(34, 8).Value = 2522.3333
$ws.Cells.Item(34, 9).Value = 2462.625
$ws.Cells.Item(34, 11).Value = 2462.625
$ws.Cells.Item(34, 13).Value = -2260.625
$ws.Cells.Item(51, 8).Value = 15466.333
$ws.Cells.Item(61, 8).Value = 15466.333
$ws.Cells.Item(74, 8).Value = 45650
$ws.Cells.Item(74, 10).Value = 45650
$ws.Cells.Item(74, 12).Value = 45650
$ws.Cells.Item(74, 14).Value = -47398
$ws.Cells.Item(77, 8).Value = 45650
$ws.Cells.Item(77, 10).Value = 45650
$ws.Cells.Item(77, 12).Value = 136950
$ws.Cells.Item(77, 14).Value = -145686
$ws.Cells.Item(86, 8).Value = 8000
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 8000
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 8000
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(86, 14).Value = -10246
$ws.Cells.Item(89, 8).Value = 8000
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 8000
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 40000
$ws.Cells.Item(89, 13).ClearContents()
$ws.Cells.Item(89, 14).Value = -51232
$ws.Cells.Item(99, 8).Value = 2618
$ws.Cells.Item(99, 9).Value = 2618
$ws.Cells.Item(99, 11).Value = 2618
$ws.Cells.Item(99, 13).Value = -1120
$ws.Cells.Item(126, 8).Value = 2618
$ws.Cells.Item(126, 9).Value = 2618
$ws.Cells.Item(126, 11).Value = 7854
$ws.Cells.Item(126, 13).Value = -5384
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 104073.766
$ws.Cells.Item(7, 9).Value = 112740.914
$ws.Cells.Item(7, 10).Value = 68
$ws.Cells.Item(7, 11).Value = 338222.742
$ws.Cells.Item(7, 12).Value = 204
$ws.Cells.Item(7, 13).Value = -338110.742
$ws.Cells.Item(7, 14).Value = -428
$ws.Cells.Item(23, 8).Value = 497.16666
$ws.Cells.Item(23, 10).Value = 693.5
$ws.Cells.Item(23, 12).Value = 2080.5
$ws.Cells.Item(23, 14).Value = -2550.5
$ws.Cells.Item(131, 8).Value = 2766.7144
$ws.Cells.Item(131, 9).Value = 2591.75
$ws.Cells.Item(131, 11).Value = 7775.25
$ws.Cells.Item(131, 13).Value = -2735.25
$ws.Cells.Item(139, 8).Value = 3732.8572
$ws.Cells.Item(139, 9).Value = 3459.4
$ws.Cells.Item(139, 11).Value = 10378.2
$ws.Cells.Item(139, 13).Value = -5238.200000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6729.4
$ws.Cells.Item(70, 9).Value = 8332.666999999999
$ws.Cells.Item(70, 10).Value = 4324.5
$ws.Cells.Item(70, 11).Value = 8332.666999999999
$ws.Cells.Item(70, 12).Value = 4324.5
$ws.Cells.Item(70, 13).Value = -8062.666999999999
$ws.Cells.Item(70, 14).Value = -4864.5
$ws.Cells.Item(73, 8).Value = 6729.4
$ws.Cells.Item(73, 9).Value = 8332.666999999999
$ws.Cells.Item(73, 10).Value = 4324.5
$ws.Cells.Item(73, 11).Value = 8332.666999999999
$ws.Cells.Item(73, 12).Value = 4324.5
$ws.Cells.Item(73, 13).Value = -7396.666999999999
$ws.Cells.Item(73, 14).Value = -6196.5
$ws.Cells.Item(122, 8).Value = 2807.6667
$ws.Cells.Item(122, 9).Value = 1758.625
$ws.Cells.Item(122, 11).Value = 5275.875
$ws.Cells.Item(122, 13).Value = -2825.875
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 700
$ws.Cells.Item(40, 9).Value = 700
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 700
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 14).ClearContents()
$ws.Cells.Item(40, 13).Value = -564
$ws.Cells.Item(46, 8).Value = 2172.5557
$ws.Cells.Item(46, 9).Value = 1214
$ws.Cells.Item(46, 11).Value = 1214
$ws.Cells.Item(46, 13).Value = -1026
$ws.Cells.Item(122, 8).Value = 4173.375
$ws.Cells.Item(122, 9).Value = 4078.6
$ws.Cells.Item(122, 11).Value = 12235.8
$ws.Cells.Item(122, 13).Value = -9785.799999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 227795.6
$ws.Cells.Item(14, 9).Value = 434989.5
$ws.Cells.Item(14, 11).Value = 434989.5
$ws.Cells.Item(14, 13).Value = -434821.5
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(45, 8).Value = 40879.332
$ws.Cells.Item(45, 9).Value = 32996.25
$ws.Cells.Item(45, 10).Value = 47185.8
$ws.Cells.Item(45, 11).Value = 32996.25
$ws.Cells.Item(45, 12).Value = 47185.8
$ws.Cells.Item(45, 13).Value = -32505.25
$ws.Cells.Item(45, 14).Value = -48167.8
$ws.Cells.Item(54, 8).Value = 60000
$ws.Cells.Item(54, 10).Value = 60000
$ws.Cells.Item(54, 12).Value = 60000
$ws.Cells.Item(54, 14).Value = -61040
$ws.Cells.Item(122, 8).Value = 2757.7646
$ws.Cells.Item(122, 9).Value = 2655.8572
$ws.Cells.Item(122, 11).Value = 7967.571599999999
$ws.Cells.Item(122, 13).Value = -5517.571599999999
$ws.Cells.Item(126, 8).Value = 2192.9412
$ws.Cells.Item(126, 9).Value = 2105.4666
$ws.Cells.Item(126, 10).Value = 2849
$ws.Cells.Item(126, 11).Value = 6316.399800000001
$ws.Cells.Item(126, 12).Value = 8547
$ws.Cells.Item(126, 13).Value = -3846.399800000001
$ws.Cells.Item(126, 14).Value = -13487
$ws.Cells.Item(135, 8).Value = 61333
$ws.Cells.Item(135, 9).Value = 54999
$ws.Cells.Item(135, 10).Value = 64500
$ws.Cells.Item(135, 11).Value = 54999
$ws.Cells.Item(135, 12).Value = 64500
$ws.Cells.Item(135, 14).Value = -74640
$ws.Cells.Item(135, 13).Value = -49929
